$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 120, shifting existing rows 120-224 down to 121-225
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new data record
$ws.Range("A120").Value = 7
$ws.Range("B120").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C120").Value = "Ñuble"
$ws.Range("D120").Value = 44589
$ws.Range("D120").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E120").Value = 16
$ws.Range("F120").Value = 100112008
$ws.Range("G120").Value = "Coliflor"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 240
$ws.Range("K120").Value = 800
$ws.Range("L120").Value = 850
$ws.Range("M120").Value = 825
$ws.Range("N120").Value = "$/unidad"
$ws.Range("O120").Value = "Provincia de Diguillín"
$ws.Range("P120").Value = 825
$ws.Range("Q120").Value = 1
$ws.Range("R120").Value = "Hortaliza"
